# Update cryptocurrency price (D) and 1h volume change (E) columns
# with latest scraped values, preserving original text formatting
# (prices/percentages are stored as text, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "96.440.04"
    "E2" = "  +0.08%  "
    "D3" = "3.676.54"
    "E3" = "  -0.28%  "
    "D4" = "2.46"
    "E4" = "  +32.21%  "
    "E5" = "  -0.05%  "
    "D6" = "228.40"
    "E6" = "  -3.04%  "
    "D7" = "651.38"
    "E7" = "  +0.38%  "
    "D8" = "0.438"
    "E8" = "  +2.59%  "
    "E9" = "  +9.98%  "
    "E10" = "  -0.04%  "
    "D11" = "3.672.85"
    "E11" = "  -0.31%  "
    "D12" = "47.88"
    "E12" = "  +8.78%  "
    "E13" = "  +2.95%  "
    "E14" = "  -4.38%  "
    "E15" = "  -1.41%  "
    "D16" = "4.363.48"
    "E16" = "  -0.25%  "
    "D17" = "96.179.46"
    "E17" = "  -0.10%  "
    "D18" = "8.89"
    "E18" = "  +1.62%  "
    "D19" = "3.677.31"
    "E19" = "  -0.34%  "
    "D20" = "19.50"
    "E20" = "  +4.94%  "
    "D21" = "12.81"
    "E21" = "  -0.75%  "
    "D22" = "0.548"
    "E22" = "  +10.41%  "
    "D23" = "530.56"
    "E23" = "  +2.84%  "
    "E24" = "  -1.48%  "
    "D25" = "0.249"
    "E25" = "  +43.84%  "
    "D26" = "121.60"
    "E26" = "  +21.11%  "
    "E27" = "  +0.21%  "
    "D28" = "6.83"
    "E28" = "  -0.85%  "
    "D29" = "3.875.82"
    "E29" = "  -0.39%  "
    "D30" = "13.07"
    "E30" = "  -0.08%  "
    "D31" = "13.34"
    "E31" = "  +10.66%  "
    "D32" = "2.99"
    "E32" = "  +0.38%  "
    "D33" = "0.999"
    "E33" = "  -0.09%  "
    "D34" = "0.185"
    "E34" = "  +1.00%  "
    "D35" = "33.10"
    "E35" = "  +3.57%  "
    "E36" = "  -1.68%  "
    "D37" = "1.00"
    "E37" = "  +0.11%  "
    "D38" = "0.610"
    "E38" = "  +4.70%  "
    "D39" = "605.85"
    "E39" = "  -6.63%  "
    "E40" = "  -0.01%  "
    "D41" = "8.40"
    "E41" = "  -3.60%  "
    "D42" = "7.14"
    "E42" = "  +4.00%  "
    "D43" = "0.508"
    "E43" = "  +19.47%  "
    "E44" = "  +2.54%  "
    "E45" = "  +11.87%  "
    "D46" = "40.26"
    "E46" = "  +0.94%  "
    "E47" = "  -2.95%  "
    "D48" = "0.963"
    "E48" = "  +1.67%  "
    "D49" = "8.99"
    "E49" = "  +6.98%  "
    "D50" = "2.29"
    "E50" = "  +1.09%  "
    "D51" = "23.52"
    "E51" = "  -0.18%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "96.440.04",
    # "1.00", "  +0.08%  ") are not auto-coerced into numbers/dates,
    # then restore the default "Normal" style so no stray formatting
    # is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
